# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates column G ("K") values for rows 2-35 with the recalculated results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(4,2,6,7,5,5,2,8,9,1,3,9,8,4,3,6,11,6,7,3,9,4,10,7,6,6,6,4,8,3,3,1,1,1)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
